$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.886.13"
$ws.Range("E2").Value = "  -4.18%  "
$ws.Range("D3").Value = "2.450.02"
$ws.Range("E3").Value = "  -3.46%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'309.68"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").Value = "'94.06"
$ws.Range("E6").Value = "  -6.77%  "
$ws.Range("D7").Value = "'0.549"
$ws.Range("E7").Value = "  -3.89%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.501"
$ws.Range("E9").Value = "  -5.07%  "
$ws.Range("D10").Value = "'33.40"
$ws.Range("E10").Value = "  -7.82%  "
$ws.Range("D11").Value = "'0.0779"
$ws.Range("E11").Value = "  -3.07%  "
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("D13").Value = "'6.93"
$ws.Range("E13").Value = "  -5.56%  "
$ws.Range("D14").Value = "2.824.88"
$ws.Range("E14").Value = "  -3.49%  "
$ws.Range("D15").Value = "2.453.68"
$ws.Range("E15").Value = "  -3.64%  "
$ws.Range("D16").Value = "'14.41"
$ws.Range("E16").Value = "  -8.64%  "
$ws.Range("D17").Value = "'0.786"
$ws.Range("E17").Value = "  -3.32%  "
$ws.Range("D18").Value = "40.928.89"
$ws.Range("E18").Value = "  -4.04%  "
$ws.Range("D19").Value = "'6.35"
$ws.Range("E19").Value = "  -6.15%  "
$ws.Range("D20").Value = "0.0₃0911"
$ws.Range("E20").Value = "  -4.27%  "
$ws.Range("D21").Value = "'11.52"
$ws.Range("E21").Value = "  -5.95%  "
$ws.Range("D22").Value = "'66.90"
$ws.Range("E22").Value = "  -3.50%  "
$ws.Range("D23").Value = "'236.76"
$ws.Range("E23").Value = "  -3.03%  "
$ws.Range("D24").Value = "'2.76"
$ws.Range("E24").Value = "  -4.61%  "
$ws.Range("D25").Value = "'1.92"
$ws.Range("E25").Value = "  -6.30%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "'24.50"
$ws.Range("E27").Value = "  -6.34%  "
$ws.Range("E28").Value = "  -3.87%  "
$ws.Range("D29").Value = "'9.66"
$ws.Range("E29").Value = "  -5.08%  "
$ws.Range("D30").Value = "'36.09"
$ws.Range("E30").Value = "  -8.30%  "
$ws.Range("D31").Value = "'152.95"
$ws.Range("E31").Value = "  -1.88%  "
$ws.Range("D32").Value = "'5.58"
$ws.Range("E32").Value = "  -3.40%  "
$ws.Range("E33").Value = "  -1.03%  "
$ws.Range("D34").Value = "'0.0749"
$ws.Range("D35").Value = "'2.51"
$ws.Range("E35").Value = "  -9.03%  "
$ws.Range("D36").Value = "'3.02"
$ws.Range("E36").Value = "  -5.46%  "
$ws.Range("D37").Value = "'17.18"
$ws.Range("E37").Value = "  -6.18%  "
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value = "'0.114"
$ws.Range("E39").Value = "  -4.46%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.102"
$ws.Range("E40").Value = "  -8.74%  "
$ws.Range("D41").Value = "'4.13"
$ws.Range("E41").Value = "  -4.24%  "
$ws.Range("D42").Value = "'21.13"
$ws.Range("E42").Value = "  -5.34%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").Value = "1.965.86"
$ws.Range("E44").Value = "  -0.74%  "
$ws.Range("E45").Value = "  -5.20%  "
$ws.Range("D46").Value = "'3.04"
$ws.Range("E46").Value = "  -8.35%  "
$ws.Range("D47").Value = "'8.71"
$ws.Range("E47").Value = "  -1.99%  "
$ws.Range("D48").Value = "'76.45"
$ws.Range("E48").Value = "  -5.46%  "
$ws.Range("D49").Value = "'97.24"
$ws.Range("E49").Value = "  -3.47%  "
$ws.Range("D50").Value = "'69.00"
$ws.Range("E50").Value = "  -5.00%  "
$ws.Range("D51").Value = "'0.179"
$ws.Range("E51").Value = "  -6.64%  "
